$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.08150863212614792
$ws.Range("D2").Value = 0.7436556217677459

$ws.Range("C3").Value = 0.1568405326489421
$ws.Range("D3").Value = 0.5906268483453735

$ws.Range("C4").Value = 0.2346779462782754
$ws.Range("D4").Value = 0.4772238367102605

$ws.Range("C5").Value = 0.3273595898062373
$ws.Range("D5").Value = 0.3701799637633356

$ws.Range("C6").Value = 0.3124981261414224
$ws.Range("D6").Value = 0.3855804851688722

$ws.Range("C7").Value = 0.2570425433549289
$ws.Range("D7").Value = 0.4490903601106034

$ws.Range("C8").Value = 0.1655040290119493
$ws.Range("D8").Value = 0.5764867229046448

$ws.Range("C9").Value = 0.08776696834527142
$ws.Range("D9").Value = 0.727376560556585

$ws.Range("C10").Value = 0.06636296460780945
$ws.Range("D10").Value = 0.7909867264656207

$ws.Range("C11").Value = 0.7245
$ws.Range("D11").Value = 0.09909999999999999
